$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "ハム太郎" post row (row 832). Excel will shift all subsequent
# rows up by one, which matches the target diff (old row 833 -> new row 832,
# ..., old row 863 -> new row 862) and the updated sheet dimension.
$ws.Rows("832").Delete()
